# semana 41 de 2024
# Update the "Esperado" (C), "Observado" (D) and "valor p" (E) columns
# for the poisson worksheet with the refreshed weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 1

# Row 3
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0.37

# Row 5
$ws.Cells.Item(5, 3).Value = 6
$ws.Cells.Item(5, 4).Value = 3
$ws.Cells.Item(5, 5).Value = 0.09

# Row 6
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(6, 4).Value = 46

# Row 7
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 0.27

# Row 11
$ws.Cells.Item(11, 3).Value = 48
$ws.Cells.Item(11, 4).Value = 32
$ws.Cells.Item(11, 5).Value = 0

# Row 12
$ws.Cells.Item(12, 3).Value = 1
$ws.Cells.Item(12, 5).Value = 0.37

# Row 13
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.14

# Row 14
$ws.Cells.Item(14, 3).Value = 5
$ws.Cells.Item(14, 4).Value = 2
$ws.Cells.Item(14, 5).Value = 0.08

# Row 15
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 1

# Row 18
$ws.Cells.Item(18, 3).Value = 10
$ws.Cells.Item(18, 5).Value = 0

# Row 19
$ws.Cells.Item(19, 3).Value = 9
$ws.Cells.Item(19, 4).Value = 3
$ws.Cells.Item(19, 5).Value = 0.01

# Row 20
$ws.Cells.Item(20, 3).Value = 2
$ws.Cells.Item(20, 5).Value = 0.14

# Row 21
$ws.Cells.Item(21, 3).Value = 6
$ws.Cells.Item(21, 5).Value = 0.01

# Row 24
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 0.14

# Row 27
$ws.Cells.Item(27, 3).Value = 5
$ws.Cells.Item(27, 4).Value = 5
$ws.Cells.Item(27, 5).Value = 0.18

# Row 28
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = 0

# Row 34
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 5).Value = 0.37

# Row 35
$ws.Cells.Item(35, 3).Value = 6
$ws.Cells.Item(35, 4).Value = 5
$ws.Cells.Item(35, 5).Value = 0.16

# Row 36
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 5).Value = 0

# Row 37
$ws.Cells.Item(37, 3).Value = 8
$ws.Cells.Item(37, 4).Value = 12
$ws.Cells.Item(37, 5).Value = 0.05
